# 🔄 Actualización automática del mapa (2025-08-11 11:26:58)
#
# The upstream data source dropped the old case "6193" (POLA 591) row,
# so the row at A81:P81 in the "AYKO" sheet needs to be removed entirely,
# shifting every following row up by one (old row 82 becomes new row 81,
# ... old row 87 becomes new row 86). This also shrinks the used range
# from A1:P87 down to A1:P86.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 81 (Caso 6193 / "POLA 591"), shifting rows 82-87
# up by one position.
$ws.Rows("81:81").Delete()
